$wb = $excel.ActiveWorkbook
$ws3 = $wb.Worksheets.Item(3)
$c = $ws3.Range("B19")
$c.NumberFormat = "[$-409]m/d/yyyy;@"
